$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the empty placeholder cells left over in row 4 (H4:L4) ---
$ws.Range("H4:L4").ClearContents()

# --- New row 5 ---
$ws.Range("A5").Value = "'202110370311392"
$ws.Range("B5").Value = "MUHAMMAD IBNU"
$ws.Range("C5").Value = "Jaringan Komputer"
$ws.Range("D5").Value = "Pemrogaman Mobile"
$ws.Range("E5").Value = "Sistem Operasi"
$ws.Range("M5").Value = 11
$ws.Range("N5").Value = "'392"

# --- New row 6 ---
$ws.Range("A6").Value = "'202110370311393"
$ws.Range("B6").Value = "NABILA AZ-ZAHRO"
$ws.Range("C6").Value = "Metode Penelitian"
$ws.Range("D6").Value = "Jaringan Komputer"
$ws.Range("E6").Value = "Basisdata"
$ws.Range("M6").Value = 8
$ws.Range("N6").Value = "'393"

# --- New row 7 ---
$ws.Range("A7").Value = "'202110370311433"
$ws.Range("B7").Value = "AL GHOZI MUHAMMAD FATUR RAHMAN"
$ws.Range("C7").Value = "Pemrogrman Website"
$ws.Range("D7").Value = "Sistem Operasi"
$ws.Range("E7").Value = "Pengantar Game"
$ws.Range("F7").Value = "'"
$ws.Range("G7").Value = "'"
$ws.Range("H7").Value = "'"
$ws.Range("I7").Value = "'"
$ws.Range("J7").Value = "'"
$ws.Range("K7").Value = "'"
$ws.Range("L7").Value = "'"
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = "'433"

# --- Update selection to mirror the diff (active cell moved to A5) ---
$ws.Range("A5").Select()
